# Append new Lancers listing scraped at 2025-09-07 01:48:15
# Inserts a new row 5 (shifting existing rows 5-9 down to 6-10),
# refreshes the "取得日時" timestamp for every data row, widens column B,
# and rewires the F-column hyperlinks so each URL cell keeps the right target.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B (35 -> 42). The engine stores column widths with a fixed
# ~0.8333 character padding offset on save, so back that off here to land
# on exactly width="42" in the saved OOXML (same reason the untouched
# columns' literal widths - e.g. 35 - look "round" in the original file).
$ws.Columns.Item(2).ColumnWidth = 41.166666666666664

# Insert a fresh row above the current row 5; this shifts rows 5-9 (and
# their cell styles) down to rows 6-10, keeping the sheet's row count growth
# in sync with the new dimension (A1:H10).
$ws.Rows.Item(5).Insert()

# New entry that now occupies row 5.
$ws.Cells.Item(5, 1).Value = "2025-09-07 01:48:15"
$ws.Cells.Item(5, 2).Value = "【急募】Instagram投稿を自動でGoogleビジネスに連携するMEOツール"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5388589"
$ws.Cells.Item(5, 7).Value = 68
$ws.Cells.Item(5, 8).Value = "◆ツール"

# Refresh the capture timestamp on every other data row (2-4, 6-10).
$ws.Cells.Item(2, 1).Value = "2025-09-07 01:48:15"
$ws.Cells.Item(3, 1).Value = "2025-09-07 01:48:15"
$ws.Cells.Item(4, 1).Value = "2025-09-07 01:48:15"
$ws.Cells.Item(6, 1).Value = "2025-09-07 01:48:15"
$ws.Cells.Item(7, 1).Value = "2025-09-07 01:48:15"
$ws.Cells.Item(8, 1).Value = "2025-09-07 01:48:15"
$ws.Cells.Item(9, 1).Value = "2025-09-07 01:48:15"
$ws.Cells.Item(10, 1).Value = "2025-09-07 01:48:15"

# The row Insert() shifted cell values/styles down correctly, but this
# engine's Hyperlinks collection does not move along with it, so rebuild
# the F-column hyperlinks from scratch in the correct order.
$ws.Hyperlinks.Delete()

$urls = @(
    "https://www.lancers.jp/work/detail/5388502",
    "https://www.lancers.jp/work/detail/5388329",
    "https://www.lancers.jp/work/detail/5388066",
    "https://www.lancers.jp/work/detail/5388589",
    "https://www.lancers.jp/work/detail/5388547",
    "https://www.lancers.jp/work/detail/5388189",
    "https://www.lancers.jp/work/detail/5385681",
    "https://www.lancers.jp/work/detail/5388228",
    "https://www.lancers.jp/work/detail/5388482"
)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $urls[$i])
}
